$d = $word.ActiveDocument

# --- Paragraph 1: "Do some exercises" ---
$p1 = $d.Paragraphs.Item(1)
$p1.LineSpacingRule = 5   # wdLineSpaceMultiple
$p1.LineSpacing = 24      # 24pt -> stored as 480 (24*20), lineRule="auto" (double spacing)
$p1.Range.Font.Size = 24
$p1.Range.Font.SizeBi = 24

# --- Paragraph 2: "Discuss about the project data" ---
$p2 = $d.Paragraphs.Item(2)
$p2.LineSpacingRule = 5
$p2.LineSpacing = 24
$p2.Range.Font.Size = 24
$p2.Range.Font.SizeBi = 24

# --- Paragraph 3: currently empty except for the "_GoBack" bookmark ---
# Insert the trailing text right after the bookmark (end), then the leading text
# right before the bookmark (start), so the bookmark stays positioned between
# the two new runs, matching the original document structure.
$bm = $d.Bookmarks.Item("_GoBack")
$bmEnd = $bm.End
$rAfter = $d.Range($bmEnd, $bmEnd)
$rAfter.InsertAfter("material that is important")

$bm2 = $d.Bookmarks.Item("_GoBack")
$bmStart = $bm2.Start
$rBefore = $d.Range($bmStart, $bmStart)
$rBefore.InsertBefore("Trying to find answers on the lectures’ ")

$p3 = $d.Paragraphs.Item(3)
$p3.LineSpacingRule = 5
$p3.LineSpacing = 24
$p3.Range.Font.Size = 24
$p3.Range.Font.SizeBi = 24
